# Applies updated organization/success-rate data to the Summary sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# Column E width 19 -> 20 (OOXML "width" units).
# The COM ColumnWidth setter adds the standard 5px padding term (5/6 in this
# engine's char-width units) before writing the raw OOXML width, so back
# out that offset to land exactly on 20.
$ws.Columns.Item(5).ColumnWidth = 20 - (5/6)

# Row data: row, Organizations_Found (C), Success_Rate_% (E), Data_Quality (F)
$updates = @(
    @{Row=2;  C=3;  E=75},
    @{Row=3;  C=20; E=25},
    @{Row=4;  C=4;  E=80},
    @{Row=5;  C=69; E=363.1578947368421},
    @{Row=6;  C=48; E=342.8571428571428},
    @{Row=7;  C=4;  E=66.66666666666666},
    @{Row=8;  C=31; E=193.75},
    @{Row=9;  C=32; E=42.66666666666667},
    @{Row=10; C=9;  E=90},
    @{Row=11; C=14; E=93.33333333333333}
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 3).Value = $u.C
    $ws.Cells.Item($r, 5).Value = $u.E
    $ws.Cells.Item($r, 6).Value = "Issues"
}

# Row 12 keeps its numeric values but Data_Quality changes to Issues
$ws.Cells.Item(12, 6).Value = "Issues"
